$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update row 2 (was Hello2/There2/Test2 content source now sits at row 2)
$ws.Range("A2").Value = "Hello2"
$ws.Range("B2").Value = "There2"
$ws.Range("C2").Value = "Test2"

# Update row 3 with the "2" variants of the former row-3 words
# (written right-to-left so the shared-string table records new
# entries in the same order Excel produced them: here2, data2, New2)
$ws.Range("C3").Value = "here2"
$ws.Range("B3").Value = "data2"
$ws.Range("A3").Value = "New2"

# Add new row 4
$ws.Range("A4").Value = "blah"
$ws.Range("B4").Value = "blah"

# Update selection to B4
$ws.Range("B4").Select()
